$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '46.702.14'
$ws.Range('E2').Value = '  +3.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.262.66'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '299.40'
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.68'
$ws.Range('E6').Value = '  +4.75%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.560'
$ws.Range('E7').Value = '  -1.26%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.508'
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.17'
$ws.Range('E10').Value = '  +2.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0806'
$ws.Range('E11').Value = '  +1.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.07'
$ws.Range('E12').Value = '  -2.15%  '
$ws.Range('E13').Value = '  -1.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.610.19'
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.262.30'
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '13.61'
$ws.Range('E16').Value = '  -0.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '46.753.72'
$ws.Range('E17').Value = '  +3.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.791'
$ws.Range('E18').Value = '  -2.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.73'
$ws.Range('E19').Value = '  -3.64%  '
$ws.Range('E20').Value = '  +3.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.80'
$ws.Range('E21').Value = '  -3.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.39'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '246.57'
$ws.Range('E23').Value = '  +2.61%  '
$ws.Range('E24').Value = '  -3.20%  '
$ws.Range('E25').Value = '  -0.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.85'
$ws.Range('E26').Value = '  -3.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '41.56'
$ws.Range('E27').Value = '  +0.89%  '
$ws.Range('E28').Value = '  -0.56%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.61'
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.41'
$ws.Range('E30').Value = '  +3.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.83'
$ws.Range('E31').Value = '  +10.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '145.96'
$ws.Range('E32').Value = '  -3.79%  '
$ws.Range('E33').Value = '  +12.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.37'
$ws.Range('E34').Value = '  -2.88%  '
$ws.Range('E35').Value = '  -3.02%  '
$ws.Range('E36').Value = '  +9.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.115'
$ws.Range('E37').Value = '  -2.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '15.60'
$ws.Range('E38').Value = '  +15.57%  '
$ws.Range('E39').Value = '  -4.19%  '
$ws.Range('E40').Value = '  -2.30%  '
$ws.Range('E41').Value = '  -4.85%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.10'
$ws.Range('E42').Value = '  -3.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.08%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.785.85'
$ws.Range('E44').Value = '  +1.08%  '
$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.63'
$ws.Range('E45').Value = '  +18.67%  '
$ws.Range('E46').Value = '  -3.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '71.18'
$ws.Range('E47').Value = '  +1.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.184'
$ws.Range('E48').Value = '  -3.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.80'
$ws.Range('E49').Value = '  +1.19%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.85'
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '94.05'
$ws.Range('E51').Value = '  -1.51%  '
